$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 638, shifting the existing rows (old 638-679) down to 639-680.
$ws.Rows.Item(638).Insert()

# Fill in the new row's data. The date/weekday columns hold plain text in this
# sheet (not real Excel dates), so prefix the date with an apostrophe to stop
# it being auto-recognised as a date serial, then reset the cell style back to
# Normal so no stray number-format style is left applied to the cell.
$ws.Range("A638").Value = "'2026/01/14"
$ws.Range("A638").Style = "Normal"
$ws.Range("B638").Value = "水"
$ws.Range("C638").Value = 20
$ws.Range("D638").Value = 201
